$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "DQ_Report": insert 6 new rows (rows 2-7) with new data,
# pushing the existing data rows down to rows 8-17.
# ---------------------------------------------------------------
$wsReport = $wb.Worksheets.Item("DQ_Report")
$wsReport.Rows("2:7").Insert()
$wsReport.Range("A2:D7").ClearFormats()

$newRows = @(
    @("P_1057020", "E70.0", 79254, "Relation  E70.0 - 79254  ist im BfArM nicht vorhanden "),
    @("P_1057020", "E70.0", 79254, "Relation  E70.0 - 79254  ist im BfArM nicht vorhanden "),
    @("P_1695115", "E70.0", 79254, "Relation  E70.0 - 79254  ist im BfArM nicht vorhanden "),
    @("P_1695115", "E70.0", 79254, "Relation  E70.0 - 79254  ist im BfArM nicht vorhanden "),
    @("P_1897170", "E70.0", 79254, "Relation  E70.0 - 79254  ist im BfArM nicht vorhanden "),
    @("P_1897170", "E70.0", 79254, "Relation  E70.0 - 79254  ist im BfArM nicht vorhanden ")
)

$r = 2
foreach ($row in $newRows) {
    $wsReport.Cells.Item($r, 1).Value = $row[0]
    $wsReport.Cells.Item($r, 2).Value = $row[1]
    $wsReport.Cells.Item($r, 3).Value = $row[2]
    $wsReport.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}

# ---------------------------------------------------------------
# Sheet "Statistik": rename some headers, add icdRd_no / pt_no
# columns, and update the uniqueness_rate value.
# ---------------------------------------------------------------
$wsStat = $wb.Worksheets.Item("Statistik")
$wsStat.Range("C1").Value = "completness_rate"
$wsStat.Range("D1").Value = "orphaCoding_completeness"
$wsStat.Range("E1").Value = "uniqueness_rate"
$wsStat.Range("F1").Value = "icdRd_no"
$wsStat.Range("G1").Value = "pt_no"

$wsStat.Range("E2").Value = 96.59999999999999
$wsStat.Range("F2").Value = 36
$wsStat.Range("G2").Value = 382

Write-Host "Done"
